# Fix typos in main block assignments
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of cell -> corrected value (typo fixes within the stimuli image columns)
$fixes = @{
    "C3"  = "stimuli/bead_y.PNG"
    "C4"  = "stimuli/bead_b.PNG"
    "D4"  = "stimuli/bead_y.PNG"
    "C5"  = "stimuli/bead_g.PNG"
    "D5"  = "stimuli/bead_b.PNG"
    "E5"  = "stimuli/bead_y.PNG"
    "C6"  = "stimuli/bead_b.PNG"
    "D6"  = "stimuli/bead_g.PNG"
    "E6"  = "stimuli/bead_b.PNG"
    "F6"  = "stimuli/bead_y.PNG"
    "D7"  = "stimuli/bead_b.PNG"
    "E7"  = "stimuli/bead_g.PNG"
    "F7"  = "stimuli/bead_b.PNG"
    "C8"  = "stimuli/bead_b.PNG"
    "E8"  = "stimuli/bead_b.PNG"
    "F8"  = "stimuli/bead_g.PNG"
    "C9"  = "stimuli/bead_g.PNG"
    "D9"  = "stimuli/bead_b.PNG"
    "F9"  = "stimuli/bead_b.PNG"
    "C10" = "stimuli/bead_b.PNG"
    "D10" = "stimuli/bead_g.PNG"
    "E10" = "stimuli/bead_b.PNG"
    "D11" = "stimuli/bead_b.PNG"
    "E11" = "stimuli/bead_g.PNG"
    "F11" = "stimuli/bead_b.PNG"
    "E12" = "stimuli/bead_b.PNG"
    "F12" = "stimuli/bead_g.PNG"
    "C13" = "stimuli/bead_b.PNG"
    "F13" = "stimuli/bead_b.PNG"
    "D14" = "stimuli/bead_b.PNG"
    "C15" = "stimuli/bead_y.PNG"
    "E15" = "stimuli/bead_b.PNG"
    "C16" = "stimuli/bead_y.PNG"
    "D16" = "stimuli/bead_y.PNG"
    "F16" = "stimuli/bead_b.PNG"
}

foreach ($addr in $fixes.Keys) {
    $ws.Range($addr).Value = $fixes[$addr]
}

# Restore view state to match the saved workbook (scroll position, zoom, selection)
$ws.Activate()
$window = $excel.ActiveWindow
$window.Zoom = 113
$ws.Range("C17:F21").Select()
$window.ScrollRow = 5
$window.ScrollColumn = 1
